$d = $word.ActiveDocument

# wdAlignParagraphCenter
$wdAlignParagraphCenter = 1

# --- 1. Center the two standalone paragraphs surrounding the table ---
# The blank paragraph right after the title (before the table)
$d.Paragraphs.Item(2).Alignment = $wdAlignParagraphCenter

# The blank paragraph right after the table (last paragraph in the body)
$d.Paragraphs.Item($d.Paragraphs.Count).Alignment = $wdAlignParagraphCenter

# --- 2. Center every data-row cell in the table (the header row already is) ---
$t = $d.Tables.Item(1)
for ($r = 2; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    for ($c = 1; $c -le $row.Cells.Count; $c++) {
        $row.Cells.Item($c).Range.Paragraphs.Item(1).Alignment = $wdAlignParagraphCenter
    }
}

# --- 3. Move the "_GoBack" bookmark from the "Nath1234" cell to the title paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Inserting a bookmark collapsed exactly at position 0 places its end marker in
# the following paragraph, so seed a throw-away character, anchor the bookmark
# right after it, then remove the character again.
$seed = $d.Range(0, 0)
$seed.InsertBefore("X")
$anchor = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $anchor)
$d.Range(0, 1).Delete()

Write-Output "done"
